$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextCell "D2" "58.689.64"
$ws.Range("E2").Value = "  -3.56%  "
Set-TextCell "D3" "2.563.20"
$ws.Range("E3").Value = "  -1.17%  "
Set-TextCell "D4" "0.999"
$ws.Range("E4").Value = "  -0.08%  "
Set-TextCell "D5" "506.14"
$ws.Range("E5").Value = "  -3.17%  "
Set-TextCell "D6" "144.92"
$ws.Range("E6").Value = "  -5.98%  "
Set-TextCell "D7" "0.998"
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("E8").Value = "  -6.58%  "
Set-TextCell "D9" "2.572.45"
$ws.Range("E9").Value = "  -1.09%  "
$ws.Range("E10").Value = "  -7.60%  "
$ws.Range("E11").Value = "  -3.00%  "
Set-TextCell "D12" "0.333"
$ws.Range("E12").Value = "  -4.03%  "
$ws.Range("E13").Value = "  -0.96%  "
Set-TextCell "D14" "3.014.46"
$ws.Range("E14").Value = "  -1.12%  "
Set-TextCell "D15" "58.709.16"
$ws.Range("E15").Value = "  -3.57%  "
Set-TextCell "D16" "20.61"
$ws.Range("E16").Value = "  -4.81%  "
$ws.Range("E17").Value = "  -4.81%  "
Set-TextCell "D18" "2.563.07"
$ws.Range("E18").Value = "  -1.48%  "
Set-TextCell "D19" "4.54"
$ws.Range("E19").Value = "  -4.62%  "
Set-TextCell "D20" "335.33"
$ws.Range("E20").Value = "  -5.23%  "
Set-TextCell "D21" "10.11"
$ws.Range("E21").Value = "  -4.37%  "
$ws.Range("E22").Value = "  -0.14%  "
Set-TextCell "D23" "5.96"
$ws.Range("E23").Value = "  -4.27%  "
Set-TextCell "D24" "59.68"
$ws.Range("E24").Value = "  -1.84%  "
Set-TextCell "D25" "0.408"
$ws.Range("E25").Value = "  -4.31%  "
$ws.Range("E26").Value = "  +0.09%  "
Set-TextCell "D27" "0.156"
$ws.Range("E27").Value = "  -6.10%  "
Set-TextCell "D28" "0.0₃0780"
$ws.Range("E28").Value = "  -7.66%  "
Set-TextCell "D29" "6.88"
$ws.Range("E29").Value = "  -6.73%  "
Set-TextCell "D30" "0.999"
$ws.Range("E30").Value = "  -0.09%  "
Set-TextCell "D31" "5.88"
$ws.Range("E31").Value = "  -7.43%  "
Set-TextCell "D32" "18.63"
$ws.Range("E32").Value = "  -3.78%  "
Set-TextCell "D33" "148.86"
$ws.Range("E33").Value = "  +0.53%  "
$ws.Range("E34").Value = "  -3.65%  "
Set-TextCell "D35" "3.88"
$ws.Range("E35").Value = "  -6.65%  "
Set-TextCell "D36" "0.905"
$ws.Range("E36").Value = "  -2.81%  "
$ws.Range("E37").Value = "  -7.62%  "
Set-TextCell "D38" "35.95"
$ws.Range("E38").Value = "  -1.39%  "
Set-TextCell "D39" "0.821"
$ws.Range("E39").Value = "  -4.73%  "
$ws.Range("B40").Value = "Filecoin"
$ws.Range("C40").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextCell "D40" "3.54"
$ws.Range("E40").Value = "  -6.94%  "
$ws.Range("B41").Value = "Stacks"
$ws.Range("C41").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextCell "D41" "1.38"
$ws.Range("E41").Value = "  -8.20%  "
Set-TextCell "D42" "283.15"
$ws.Range("E42").Value = "  -1.75%  "
Set-TextCell "D43" "1.00"
$ws.Range("E43").Value = "  +0.24%  "
$ws.Range("B44").Value = "Mantle"
$ws.Range("C44").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextCell "D44" "0.608"
$ws.Range("E44").Value = "  -2.07%  "
$ws.Range("B45").Value = "Stellar"
$ws.Range("C45").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextCell "D45" "0.0981"
$ws.Range("E45").Value = "  -3.28%  "
Set-TextCell "D46" "0.0532"
$ws.Range("E46").Value = "  -4.97%  "
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextCell "D47" "18.69"
$ws.Range("E47").Value = "  -4.55%  "
$ws.Range("B48").Value = "WhiteBITCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
Set-TextCell "D48" "10.32"
$ws.Range("E48").Value = "  -0.04%  "
$ws.Range("E49").Value = "  -4.42%  "
Set-TextCell "D50" "4.60"
$ws.Range("E50").Value = "  -5.75%  "
Set-TextCell "D51" "1.911.90"
$ws.Range("E51").Value = "  -2.31%  "
